$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the now-unneeded rows (5-9); this also shrinks the used range/dimension
$ws.Range("A5:F9").Delete()

# Row 2 (abc / MTB_WCS_MSE7_MS25 / Vendor / 7/16/2025 / 7/16/2025 / 1) is unchanged

# Row 3: replace with the John Vendor billing row
$ws.Range("A3").Value = "John Vendor"
$ws.Range("B3").Value = "MTB_WCS_MSE7_MS1"
$ws.Range("C3").Value = "Vendor"

# Dates are entered as plain text (matching the template's existing text-dates),
# so force the Text number format before typing, then drop back to the
# worksheet's default style afterwards so no stray formatting is left behind.
$ws.Range("D3:E3").NumberFormat = "@"
$ws.Range("D3").Value = "8/5/2025"
$ws.Range("E3").Value = "8/5/2025"
$ws.Range("D3:E3").Style = "Normal"

$ws.Range("F3").Value = 1

# Row 4: replace with the Jane Vendor billing row
$ws.Range("A4").Value = "Jane Vendor"
$ws.Range("B4").Value = "MTB_WCS_MSE7_MS2"
$ws.Range("C4").Value = "Vendor"

$ws.Range("D4:E4").NumberFormat = "@"
$ws.Range("D4").Value = "8/10/2025"
$ws.Range("E4").Value = "8/12/2025"
$ws.Range("D4:E4").Style = "Normal"

$ws.Range("F4").Value = 3
